# Scheduled market-data refresh: update cached Market Board pricing
# (columns H:N — currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# across the ALC/ARM/BSM/CRP/GSM/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 617.8333
$ws.Range("J98").Value = 2250
$ws.Range("L98").Value = 2250
$ws.Range("N98").Value = -5246

# Row 99
$ws.Range("H99").Value = 287.33334
$ws.Range("I99").Value = 244.8
$ws.Range("K99").Value = 734.4000000000001
$ws.Range("M99").Value = 763.5999999999999

# Row 111
$ws.Range("H111").Value = 166672100
$ws.Range("I111").Value = 200006350
$ws.Range("J111").Value = 800
$ws.Range("K111").Value = 600019050
$ws.Range("L111").Value = 2400
$ws.Range("M111").Value = -600015983
$ws.Range("N111").Value = -8534

# Row 116
$ws.Range("H116").Value = 1975453.1
$ws.Range("I116").Value = 7694557.5
$ws.Range("J116").Value = 3348.1035
$ws.Range("K116").Value = 7694557.5
$ws.Range("L116").Value = 3348.1035
$ws.Range("M116").Value = -7691115.5
$ws.Range("N116").Value = -10232.1035

# Row 122
$ws.Range("H122").Value = 617.8333
$ws.Range("J122").Value = 2250
$ws.Range("L122").Value = 6750
$ws.Range("N122").Value = -11650

# Row 132
$ws.Range("H132").Value = 2638.5518
$ws.Range("I132").Value = 2784.261
$ws.Range("J132").Value = 2080
$ws.Range("K132").Value = 8352.782999999999
$ws.Range("L132").Value = 6240
$ws.Range("M132").Value = -5822.782999999999
$ws.Range("N132").Value = -11300

# Row 138
$ws.Range("H138").Value = 3347.8438
$ws.Range("I138").Value = 1457.3235
$ws.Range("J138").Value = 4384.5806
$ws.Range("K138").Value = 4371.970499999999
$ws.Range("L138").Value = 13153.7418
$ws.Range("M138").Value = 768.0295000000006
$ws.Range("N138").Value = -23433.7418


$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2023.8182
$ws.Range("I2").Value = 2143.2222
$ws.Range("J2").Value = 1486.5
$ws.Range("K2").Value = 2143.2222
$ws.Range("L2").Value = 1486.5
$ws.Range("M2").Value = -2030.2222
$ws.Range("N2").Value = -1712.5

# Row 61
$ws.Range("H61").Value = 2101.7693
$ws.Range("I61").Value = 1574.8182
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 1574.8182
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -1362.8182
$ws.Range("N61").Value = -5424

# Row 74
$ws.Range("H74").Value = 1871.5714
$ws.Range("I74").Value = 2106.8462
$ws.Range("J74").Value = 1489.25
$ws.Range("K74").Value = 2106.8462
$ws.Range("L74").Value = 1489.25
$ws.Range("M74").Value = -1232.8462
$ws.Range("N74").Value = -3237.25

# Row 77
$ws.Range("H77").Value = 1871.5714
$ws.Range("I77").Value = 2106.8462
$ws.Range("J77").Value = 1489.25
$ws.Range("K77").Value = 10534.231
$ws.Range("L77").Value = 7446.25
$ws.Range("M77").Value = -6166.231
$ws.Range("N77").Value = -16182.25

# Row 116
$ws.Range("H116").Value = 2023.8182
$ws.Range("I116").Value = 2143.2222
$ws.Range("J116").Value = 1486.5
$ws.Range("K116").Value = 2143.2222
$ws.Range("L116").Value = 1486.5
$ws.Range("M116").Value = 150.7777999999998
$ws.Range("N116").Value = -6074.5

# Row 136
$ws.Range("H136").Value = 2101.7693
$ws.Range("I136").Value = 1574.8182
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 4724.4546
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -2174.4546
$ws.Range("N136").Value = -20100


$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2023.8182
$ws.Range("I3").Value = 2143.2222
$ws.Range("J3").Value = 1486.5
$ws.Range("K3").Value = 2143.2222
$ws.Range("L3").Value = 1486.5
$ws.Range("M3").Value = -2029.2222
$ws.Range("N3").Value = -1714.5

# Row 97
$ws.Range("H97").Value = 4500
$ws.Range("I97").Value = 4500
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 4500
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -3509
$ws.Range("N97").ClearContents()


$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 3156.3333
$ws.Range("I16").Value = 3002.2
$ws.Range("K16").Value = 3002.2
$ws.Range("M16").Value = -2715.2

# Row 22
$ws.Range("H22").Value = 268.84616
$ws.Range("I22").Value = 268.84616
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 268.84616
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 81.15384
$ws.Range("N22").ClearContents()

# Row 31
$ws.Range("H31").Value = 1966.6578
$ws.Range("I31").Value = 1458.2258
$ws.Range("J31").Value = 4218.2856
$ws.Range("K31").Value = 1458.2258
$ws.Range("L31").Value = 4218.2856
$ws.Range("M31").Value = -1163.2258
$ws.Range("N31").Value = -4808.2856

# Row 34
$ws.Range("H34").Value = 1966.6578
$ws.Range("I34").Value = 1458.2258
$ws.Range("J34").Value = 4218.2856
$ws.Range("K34").Value = 1458.2258
$ws.Range("L34").Value = 4218.2856
$ws.Range("M34").Value = -1256.2258
$ws.Range("N34").Value = -4622.2856

# Row 113
$ws.Range("H113").Value = 3156.3333
$ws.Range("I113").Value = 3002.2
$ws.Range("K113").Value = 3002.2
$ws.Range("M113").Value = -832.1999999999998

# Row 132
$ws.Range("H132").Value = 3225.6667
$ws.Range("I132").Value = 2122.3572
$ws.Range("J132").Value = 5432.2856
$ws.Range("K132").Value = 6367.071599999999
$ws.Range("L132").Value = 16296.8568
$ws.Range("M132").Value = -3837.071599999999
$ws.Range("N132").Value = -21356.8568


$ws = $wb.Worksheets.Item("GSM")
# Row 118
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

# Row 126
$ws.Range("H126").Value = 4766171.5
$ws.Range("J126").Value = 4442.75
$ws.Range("L126").Value = 13328.25
$ws.Range("N126").Value = -18268.25

# Row 132
$ws.Range("H132").Value = 3570.6
$ws.Range("I132").Value = 3704
$ws.Range("J132").Value = 3370.5
$ws.Range("K132").Value = 11112
$ws.Range("L132").Value = 10111.5
$ws.Range("M132").Value = -8582
$ws.Range("N132").Value = -15171.5


$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1538.5
$ws.Range("I132").Value = 1115.1904
$ws.Range("J132").Value = 2808.4285
$ws.Range("K132").Value = 3345.5712
$ws.Range("L132").Value = 8425.2855
$ws.Range("M132").Value = -815.5711999999999
$ws.Range("N132").Value = -13485.2855

# Row 136
$ws.Range("H136").Value = 7894.1875
$ws.Range("I136").Value = 9176
$ws.Range("J136").Value = 4048.75
$ws.Range("K136").Value = 27528
$ws.Range("L136").Value = 12146.25
$ws.Range("M136").Value = -24978
$ws.Range("N136").Value = -17246.25
